$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell with plain (unstyled) text style, used to restore style
# after forcing NumberFormat to text for numeric-looking values.
$textStyle = $ws.Range("D2").Style

$ws.Range("D2").Value = "61.482.81"
$ws.Range("E2").Value = "  -2.22%  "
$ws.Range("D3").Value = "3.379.73"
$ws.Range("E3").Value = "  -0.75%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.67"
$ws.Range("D5").Style = $textStyle
$ws.Range("E5").Value = "  +0.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.25"
$ws.Range("D6").Style = $textStyle
$ws.Range("E6").Value = "  +7.40%  "
$ws.Range("D8").Value = "3.380.74"
$ws.Range("E8").Value = "  -0.75%  "
$ws.Range("E9").Value = "  +0.47%  "
$ws.Range("E10").Value = "  +2.60%  "
$ws.Range("E11").Value = "  +1.75%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.388"
$ws.Range("D12").Style = $textStyle
$ws.Range("E12").Value = "  +1.98%  "
$ws.Range("E13").Value = "  -0.92%  "
$ws.Range("E14").Value = "  +1.09%  "
$ws.Range("E15").Value = "  +1.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.37"
$ws.Range("D17").Style = $textStyle
$ws.Range("E17").Value = "  +1.46%  "
$ws.Range("D18").Value = "61.548.55"
$ws.Range("E18").Value = "  -2.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.08"
$ws.Range("D19").Style = $textStyle
$ws.Range("E19").Value = "  +6.56%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.82"
$ws.Range("D20").Style = $textStyle
$ws.Range("E20").Value = "  +2.15%  "
$ws.Range("E21").Value = "  -2.40%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "385.27"
$ws.Range("D22").Style = $textStyle
$ws.Range("E22").Value = "  +1.89%  "
$ws.Range("E23").Value = "  +1.56%  "
$ws.Range("D24").Value = "3.513.45"
$ws.Range("E25").Value = "  +0.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "70.95"
$ws.Range("D26").Style = $textStyle
$ws.Range("E26").Value = "  -2.20%  "
$ws.Range("E27").Value = "  +9.97%  "
$ws.Range("E28").Value = "  +20.59%  "
$ws.Range("E29").Value = "  +11.96%  "
$ws.Range("E30").Value = "  -0.26%  "
$ws.Range("E31").Value = "  +3.46%  "
$ws.Range("E32").Value = "  +0.62%  "
$ws.Range("E33").Value = "  +3.13%  "
$ws.Range("E35").Value = "  -0.60%  "
$ws.Range("E36").Value = "  +2.55%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.60"
$ws.Range("D37").Style = $textStyle
$ws.Range("E37").Value = "  +4.88%  "
$ws.Range("E38").Value = "  +3.29%  "
$ws.Range("E39").Value = "  +3.68%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "163.13"
$ws.Range("D40").Style = $textStyle
$ws.Range("E40").Value = "  -0.76%  "
$ws.Range("E41").Value = "  +2.96%  "
$ws.Range("E42").Value = "  +13.63%  "
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("E44").Value = "  +3.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.74"
$ws.Range("D45").Style = $textStyle
$ws.Range("E45").Value = "  +0.18%  "
$ws.Range("E46").Value = "  -2.59%  "
$ws.Range("E47").Value = "  +2.93%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.56"
$ws.Range("D48").Style = $textStyle
$ws.Range("E48").Value = "  +2.54%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.94"
$ws.Range("D49").Style = $textStyle
$ws.Range("E49").Value = "  +3.81%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.42"
$ws.Range("D50").Style = $textStyle
$ws.Range("E50").Value = "  +15.16%  "
$ws.Range("E51").Value = "  +5.24%  "
